# "Ad label condition edited": add a new "Акция" (Promotion) column (F) to the
# wine list sheet, flagging two rows ("Да" / Yes) that are currently on
# promotion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- F1: new header cell, matching the style of the other header cells (A1:D1) ---
$ws.Range("F1").Value = "Акция"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)   # xlPasteFormats

# --- F4 / F9: mark these two wines as being on promotion ---
$ws.Range("F4").Value = "Да"
$ws.Range("F4").Font.Name = "Arial"
$ws.Range("F4").Font.Size = 10
$ws.Range("F4").Font.Bold = $false
$ws.Range("F4").Font.ThemeColor = 1

$ws.Range("F9").Value = "Да"
$ws.Range("F4").Copy()
$ws.Range("F9").PasteSpecial(-4122)   # xlPasteFormats, copy F4's new style onto F9

$excel.CutCopyMode = 0

# --- selection ends up on F9, like in the edited workbook ---
$ws.Range("F9").Select()
